$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ----------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("E8").Value  = 129.64
$ws1.Range("G8").Value  = 81.48
$ws1.Range("M13").Value = 2156.54
$ws1.Range("Q18").Value = 1935.16

$ws1.Range("E22").Value = "1 de 20"
$ws1.Range("G22").Value = "2 de 20"
$ws1.Range("M22").Value = "5 de 20"

# ----------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ----------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F8").Value  = 211.12
$ws2.Range("F13").Value = 2156.54
$ws2.Range("F18").Value = 1935.16
$ws2.Range("F22").Value = 14288.99

# ----------------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ----------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# FREGADEROS DE COCINA (row 4)
$ws3.Range("D4").Value = 129.64
$ws3.Range("E4").Value = 384.1910466593361
$ws3.Range("F4").Value = 0.2523008308720391

# GRIFERIAS (row 6)
$ws3.Range("D6").Value = 125.02
$ws3.Range("E6").Value = -18.2
$ws3.Range("F6").Value = 1.170380078636959

# PANELES PVC (row 14)
$ws3.Range("D14").Value = 2366.8
$ws3.Range("E14").Value = -1883.8
$ws3.Range("F14").Value = 4.900207039337475

# PORCELANATO (row 16)
$ws3.Range("D16").Value = 5620.63
$ws3.Range("E16").Value = 23911.81
$ws3.Range("F16").Value = 0.1903205424272427

# TOTAL (row 19)
$ws3.Range("D19").Value = 14288.99
$ws3.Range("E19").Value = 36098.20762291769
$ws3.Range("F19").Value = 0.2835837409917975
